# Updates cryptos list data (prices and 1h volume %) per the source diff.
# All target cells are plain-text cells (inlineStr, no explicit style) in the
# original workbook, so we force text formatting ("@") before assigning any
# numeric-looking string (to stop Excel auto-converting "174.70" -> 174.7,
# or "7.95" -> a float), then reset the style to "Normal" so no stray cell
# style attribute is introduced (matches the unstyled source cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.598.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('E2').Style = 'Normal'
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.563.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E3').Style = 'Normal'
# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E4').Style = 'Normal'
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E5').Style = 'Normal'
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E6').Style = 'Normal'
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.560.90'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('E7').Style = 'Normal'
# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E8').Style = 'Normal'
# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.79%  '
$ws.Range('E9').Style = 'Normal'
# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('E10').Style = 'Normal'
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('E11').Style = 'Normal'
# Row 12
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('E12').Style = 'Normal'
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.168.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.82%  '
$ws.Range('E13').Style = 'Normal'
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000208'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('E14').Style = 'Normal'
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('E15').Style = 'Normal'
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.566.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('E16').Style = 'Normal'
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.597.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('E17').Style = 'Normal'
# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('E18').Style = 'Normal'
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.61%  '
$ws.Range('E19').Style = 'Normal'
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('E20').Style = 'Normal'
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('E21').Style = 'Normal'
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '432.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('E22').Style = 'Normal'
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.611'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('E23').Style = 'Normal'
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('E24').Style = 'Normal'
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.705.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E25').Style = 'Normal'
# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E26').Style = 'Normal'
# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E27').Style = 'Normal'
# Row 28
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('E28').Style = 'Normal'
# Row 29
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('E29').Style = 'Normal'
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E30').Style = 'Normal'
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E31').Style = 'Normal'
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('E32').Style = 'Normal'
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.558.86'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('E33').Style = 'Normal'
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.39'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E34').Style = 'Normal'
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.154'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('E35').Style = 'Normal'
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.88'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E36').Style = 'Normal'
# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('E38').Style = 'Normal'
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.62'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('E39').Style = 'Normal'
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '174.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('E40').Style = 'Normal'
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0851'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('E41').Style = 'Normal'
# Row 42
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('E42').Style = 'Normal'
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.890'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('E43').Style = 'Normal'
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('E44').Style = 'Normal'
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '46.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.35%  '
$ws.Range('E45').Style = 'Normal'
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('E46').Style = 'Normal'
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.53'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.16%  '
$ws.Range('E47').Style = 'Normal'
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.55%  '
$ws.Range('E48').Style = 'Normal'
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('E49').Style = 'Normal'
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.21%  '
$ws.Range('E50').Style = 'Normal'
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.33%  '
$ws.Range('E51').Style = 'Normal'
